$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '63.116.23'
Set-TextValue "E2" '  +6.54%  '
Set-TextValue "D3" '3.113.45'
Set-TextValue "E3" '  +4.06%  '
Set-TextValue "E4" '  +0.09%  '
Set-TextValue "D5" '587.40'
Set-TextValue "E5" '  +4.81%  '
Set-TextValue "D6" '144.00'
Set-TextValue "E6" '  +4.38%  '
Set-TextValue "E7" '  +0.00%  '
Set-TextValue "D8" '3.100.56'
Set-TextValue "E8" '  +4.06%  '
Set-TextValue "E9" '  +2.45%  '
Set-TextValue "D10" '0.146'
Set-TextValue "E10" '  +10.66%  '
Set-TextValue "E11" '  +10.19%  '
Set-TextValue "D12" '0.468'
Set-TextValue "E12" '  +2.15%  '
Set-TextValue "D13" '0.0000245'
Set-TextValue "E13" '  +6.90%  '
Set-TextValue "D14" '35.70'
Set-TextValue "E14" '  +6.23%  '
Set-TextValue "E15" '  +0.93%  '
Set-TextValue "D16" '3.629.06'
Set-TextValue "D17" '7.27'
Set-TextValue "E17" '  -0.35%  '
Set-TextValue "D18" '63.070.78'
Set-TextValue "E18" '  +6.47%  '
Set-TextValue "D19" '3.109.20'
Set-TextValue "E19" '  +4.00%  '
Set-TextValue "D20" '453.96'
Set-TextValue "E20" '  +5.74%  '
Set-TextValue "D21" '14.11'
Set-TextValue "E21" '  +3.84%  '
Set-TextValue "D22" '0.733'
Set-TextValue "E22" '  +1.91%  '
Set-TextValue "D23" '7.60'
Set-TextValue "E23" '  +6.77%  '
Set-TextValue "D24" '13.63'
Set-TextValue "E24" '  +0.69%  '
Set-TextValue "D25" '82.26'
Set-TextValue "E25" '  +2.34%  '
Set-TextValue "E26" '  +0.15%  '
Set-TextValue "E27" '  +1.87%  '
Set-TextValue "D28" '2.71'
Set-TextValue "E28" '  +6.61%  '
Set-TextValue "B29" 'FirstDigitalUSD'
Set-TextValue "C29" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D29" '1.00'
Set-TextValue "E29" '  +0.13%  '
Set-TextValue "B30" 'RenderToken'
Set-TextValue "C30" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D30" '8.23'
Set-TextValue "E30" '  +5.06%  '
Set-TextValue "D31" '6.86'
Set-TextValue "E31" '  +12.51%  '
Set-TextValue "E32" '  +10.89%  '
Set-TextValue "D33" '27.03'
Set-TextValue "E33" '  +5.23%  '
Set-TextValue "D34" '2.38'
Set-TextValue "E34" '  +13.62%  '
Set-TextValue "D35" '0.0₃0813'
Set-TextValue "E35" '  +7.39%  '
Set-TextValue "E36" '  +4.31%  '
Set-TextValue "D37" '6.06'
Set-TextValue "E37" '  +2.03%  '
Set-TextValue "D38" '3.10'
Set-TextValue "E38" '  +13.02%  '
Set-TextValue "D39" '51.10'
Set-TextValue "E39" '  +4.63%  '
Set-TextValue "D40" '8.77'
Set-TextValue "E40" '  +1.18%  '
Set-TextValue "D41" '428.18'
Set-TextValue "E41" '  +5.57%  '
Set-TextValue "D42" '2.970.63'
Set-TextValue "E42" '  +7.15%  '
Set-TextValue "D43" '0.0373'
Set-TextValue "E43" '  +5.74%  '
Set-TextValue "D44" '0.112'
Set-TextValue "E44" '  +4.82%  '
Set-TextValue "E45" '  +9.66%  '
Set-TextValue "D46" '2.17'
Set-TextValue "E46" '  +8.33%  '
Set-TextValue "D47" '124.85'
Set-TextValue "E47" '  +1.24%  '
Set-TextValue "E48" '  -0.01%  '
Set-TextValue "D49" '34.64'
Set-TextValue "E49" '  +0.95%  '
Set-TextValue "E50" '  +1.24%  '
Set-TextValue "E51" '  +6.48%  '
